$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row for the "Reactivate incident" use case above row 26
$ws.Rows.Item(26).Insert()

$ws.Range("A26").Value = "Reactivate incident"
$ws.Range("B26").Value = 1
$ws.Range("E26").Formula = "=B26*5+C26*10+D26*15"
$ws.Range("F26").Formula = "=E26/E38"
$ws.Range("H26").Value = 1
$ws.Range("I26").Formula = "=E26*H26"

$ws.Application.Goto($ws.Range("A19"))
$ws.Range("I26").Select()
